# Slide 1 ("By Marin Fotache & Co.") title/credit box: split the single
# run into two runs, turning "Co." into "Co. " (trailing space) as its own
# run, same as the author's re-type of the tail of that line in PowerPoint.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Locate the shape that holds the "By Marin Fotache & Co." credit line
# instead of hard-coding a shape index, so the script is resilient to
# minor shape-ordering differences.
$target = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    if ($sh.HasTextFrame) {
        if ($sh.TextFrame.HasText) {
            if ($sh.TextFrame.TextRange.Text -like "*Marin Fotache*") {
                $target = $sh
            }
        }
    }
}

$tr = $target.TextFrame.TextRange

# "By Marin Fotache & Co." -> characters 20-22 are "Co."; replacing just
# that sub-range splits the original single run into two runs and inserts
# the trailing space that turns "Co." into "Co. ".
$tr.Characters(20, 3).Text = "Co. "
